# Apply crypto price/volume updates to Sheet1 (cryptos.xlsx)
# Commit: "Updated cryptos list on Fri Aug 16 11:27:09 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.328.31"
$ws.Range("E2").Value = "  -0.18%  "

# Row 3
$ws.Range("D3").Value = "2.599.37"
$ws.Range("E3").Value = "  -0.72%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.71"
$ws.Range("E5").Value = "  +0.09%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.73"
$ws.Range("E6").Value = "  +0.80%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.25%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").Value = "2.617.19"
$ws.Range("E9").Value = "  -0.41%  "

# Row 10
$ws.Range("E10").Value = "  -1.32%  "

# Row 11
$ws.Range("E11").Value = "  -1.23%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.342"
$ws.Range("E12").Value = "  +2.12%  "

# Row 13
$ws.Range("E13").Value = "  +0.00%  "

# Row 14
$ws.Range("D14").Value = "3.052.71"
$ws.Range("E14").Value = "  -0.83%  "

# Row 15
$ws.Range("D15").Value = "58.273.14"
$ws.Range("E15").Value = "  -0.20%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.41"
$ws.Range("E16").Value = "  -2.62%  "

# Row 17
$ws.Range("D17").Value = "2.639.25"
$ws.Range("E17").Value = "  +1.73%  "

# Row 18
$ws.Range("E18").Value = "  -1.30%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.45"
$ws.Range("E19").Value = "  +0.82%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.33"
$ws.Range("E20").Value = "  -1.62%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.24"
$ws.Range("E21").Value = "  -1.39%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.45"
$ws.Range("E22").Value = "  +2.87%  "

# Row 23
$ws.Range("E23").Value = "  -0.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.46"
$ws.Range("E24").Value = "  +0.62%  "

# Row 25
$ws.Range("E25").Value = "  +1.36%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.405"
$ws.Range("E26").Value = "  -2.19%  "

# Row 27
$ws.Range("D27").Value = "2.721.41"
$ws.Range("E27").Value = "  -0.45%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.06%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.06"
$ws.Range("E29").Value = "  -0.83%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0744"
$ws.Range("E30").Value = "  -5.90%  "

# Row 31
$ws.Range("E31").Value = "  -0.11%  "

# Row 32
$ws.Range("E32").Value = "  -6.28%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.59"
$ws.Range("E33").Value = "  -0.38%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.79"
$ws.Range("E34").Value = "  +0.18%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.80"
$ws.Range("E35").Value = "  -0.27%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.01"
$ws.Range("E36").Value = "  -2.09%  "

# Row 37
$ws.Range("E37").Value = "  -4.18%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.856"
$ws.Range("E38").Value = "  -4.50%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.861"
$ws.Range("E39").Value = "  +1.25%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.47"
$ws.Range("E40").Value = "  +2.42%  "

# Row 41
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.05"
$ws.Range("E41").Value = "  -0.53%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.54"
$ws.Range("E42").Value = "  -2.22%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.995"
$ws.Range("E43").Value = "  -0.33%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.608"
$ws.Range("E44").Value = "  +1.04%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "271.69"
$ws.Range("E45").Value = "  +1.28%  "

# Row 46
$ws.Range("E46").Value = "  +0.20%  "

# Row 47
$ws.Range("E47").Value = "  -1.76%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.72"
$ws.Range("E48").Value = "  -1.78%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0522"
$ws.Range("E49").Value = "  -1.65%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.69"
$ws.Range("E50").Value = "  +2.81%  "

# Row 51
$ws.Range("D51").Value = "1.968.70"
$ws.Range("E51").Value = "  -2.94%  "
